$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws1.Range("M4").Value = 97.81999999999999
$ws1.Range("M14").Value = 249.65
$ws1.Range("M32").Value = "6 de 30"

$ws2.Range("F4").Value = 1048.22
$ws2.Range("F14").Value = 249.65
$ws2.Range("F32").Value = 4537.95

$ws3.Range("D15").Value = 1835.12
$ws3.Range("E15").Value = 21623.7
$ws3.Range("F15").Value = 0.07822729361493885

$ws3.Range("D18").Value = 4527.67
$ws3.Range("E18").Value = 29407.04607548726
$ws3.Range("F18").Value = 0.1334229521746481
